$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L2").Value = "*maa://24633 (56.05), *maa://30515 (69.31), *maa://34787 (72.46), ***maa://20792 (11.93), maa://39402 (88.37), ***maa://29083 (27.78)"
$ws.Range("T2").Value = "maa://22742 (91.67), *maa://20791 (63.38)"
$ws.Range("AB2").Value = "maa://21246 (91.32), maa://36684 (95.74), ***maa://22731 (6.67)"
$ws.Range("AF2").Value = "maa://25251 (91.75), ***maa://21730 (21.74), ***maa://39501 (18.18), *maa://36675 (60.0)"
$ws.Range("H3").Value = "maa://21247 (98.44), *maa://22748 (60.0)"
$ws.Range("L3").Value = "*maa://22880 (65.95), maa://20276 (85.0), *maa://22749 (72.73)"
$ws.Range("P3").Value = "maa://21249 (94.22), maa://26254 (95.83)"
$ws.Range("X3").Value = "maa://27396 (84.52), maa://27484 (96.15), maa://27480 (82.86)"
$ws.Range("AB3").Value = "maa://24390 (93.55)"
$ws.Range("D4").Value = "maa://24632 (93.33), **maa://24303 (33.33), maa://22499 (86.67), maa://22746 (100.0)"
$ws.Range("T4").Value = "maa://32509 (98.04), maa://27295 (83.87), maa://22754 (91.67), *maa://21746 (55.81), *maa://31008 (78.05)"
$ws.Range("X4").Value = "**maa://32495 (47.91), ***maa://31785 (22.22), ***maa://36683 (28.26), maa://43217 (90.32)"
$ws.Range("D5").Value = "maa://21245 (83.49), maa://22744 (84.0)"
$ws.Range("L5").Value = "*maa://22757 (78.79)"
$ws.Range("AB5").Value = "*maa://29863 (67.65), ***maa://22752 (12.5), **maa://26013 (37.5)"
$ws.Range("D6").Value = "maa://42407 (94.59)"
$ws.Range("L6").Value = "maa://24839 (99.29)"
$ws.Range("H7").Value = "*maa://22763 (68.97)"
$ws.Range("T7").Value = "maa://21291 (84.09)"
$ws.Range("AF7").Value = "*maa://26191 (67.95), *maa://36671 (69.39), *maa://42530 (64.29)"
$ws.Range("A8").Value = "更新日期：2024.12.29 14:56:16"
$ws.Range("P8").Value = "maa://32931 (84.62), *maa://21916 (61.29), maa://23252 (92.42), maa://37496 (96.43), **maa://22759 (45.45)"
$ws.Range("P9").Value = "maa://22736 (81.52)"
$ws.Range("AB9").Value = "maa://28711 (87.25), ***maa://22740 (5.77), **maa://39938 (48.0), **maa://27377 (42.86), ***maa://25174 (19.05), maa://40166 (93.33)"
$ws.Range("AF9").Value = "maa://26206 (89.32), *maa://22865 (51.92)"
$ws.Range("L10").Value = "**maa://24395 (40.74)"
$ws.Range("T10").Value = "maa://27395 (96.07), maa://22755 (87.61), **maa://22756 (40.91), ***maa://21737 (10.61)"
$ws.Range("X10").Value = "maa://22301 (97.63), maa://22726 (100.0)"
$ws.Range("AF10").Value = "*maa://25021 (54.22), *maa://22733 (59.38), maa://22761 (100.0)"
$ws.Range("T11").Value = "maa://22747 (92.81), maa://22501 (98.59)"
$ws.Range("X11").Value = "maa://36713 (98.17)"
$ws.Range("AB11").Value = "maa://29912 (100.0), maa://22516 (88.37), *maa://20794 (52.24)"
$ws.Range("H12").Value = "maa://21867 (89.7)"
$ws.Range("X12").Value = "maa://22753 (91.46), *maa://21485 (76.87), maa://37962 (89.66)"
$ws.Range("AB12").Value = "maa://23669 (95.41), maa://36677 (92.16), maa://39872 (90.0)"
$ws.Range("D13").Value = "maa://24999 (91.79), maa://36673 (92.65), maa://25001 (85.51)"
$ws.Range("P13").Value = "maa://22676 (91.82), *maa://22583 (75.0), *maa://22500 (57.78)"
$ws.Range("L14").Value = "maa://26245 (96.53), maa://21288 (96.3), maa://39841 (95.12), maa://36682 (97.37)"
$ws.Range("D15").Value = "*maa://22743 (77.32), maa://22734 (84.03), *maa://30808 (65.08), **maa://36048 (33.33)"
$ws.Range("H15").Value = "maa://24304 (88.18), maa://21478 (91.43)"
$ws.Range("P15").Value = "maa://24762 (90.32), *maa://22727 (70.0)"
$ws.Range("D16").Value = "maa://21441 (96.35), maa://36679 (93.02), maa://37650 (96.97)"
$ws.Range("AF16").Value = "*maa://23911 (65.69), maa://27755 (92.68)"
$ws.Range("H17").Value = "maa://22430 (88.71), maa://39599 (86.84)"
$ws.Range("P17").Value = "maa://23890 (81.0), *maa://24940 (67.86)"
$ws.Range("H18").Value = "maa://24421 (90.46)"
$ws.Range("L18").Value = "maa://22466 (88.89), *maa://22732 (50.6)"
$ws.Range("AF18").Value = "*maa://24313 (57.86), **maa://29784 (44.44)"
$ws.Range("AB19").Value = "*maa://30709 (63.5), *maa://36668 (55.84)"
$ws.Range("AF19").Value = "*maa://21663 (61.19)"
$ws.Range("D20").Value = "maa://21432 (89.86), maa://25198 (93.14), *maa://20795 (51.18), maa://36680 (96.67)"
$ws.Range("L20").Value = "maa://41331 (84.69)"
$ws.Range("H21").Value = "maa://24372 (96.77)"
$ws.Range("D23").Value = "***maa://28036 (27.54), *maa://41753 (53.85)"
$ws.Range("L23").Value = "maa://39756 (94.21), maa://39875 (93.75)"
$ws.Range("X23").Value = "*maa://28503 (66.18)"
$ws.Range("D24").Value = "*maa://24368 (79.55)"
$ws.Range("X24").Value = "maa://29988 (86.84), maa://23504 (93.33), **maa://22892 (39.58), *maa://25141 (76.98), *maa://36663 (79.41), ***maa://22815 (23.08)"
$ws.Range("AF24").Value = "maa://22523 (85.64), maa://36672 (80.77), maa://29910 (92.59), **maa://21440 (34.55)"
$ws.Range("D25").Value = "maa://29753 (94.88)"
$ws.Range("H25").Value = "*maa://29063 (73.65), *maa://25311 (73.27), ***maa://22725 (4.84)"
$ws.Range("X28").Value = "maa://39929 (89.91), ***maa://39723 (14.29), maa://41749 (91.07)"
$ws.Range("AF28").Value = "maa://36660 (92.97), *maa://36701 (64.29)"
$ws.Range("L29").Value = "maa://28432 (92.9), *maa://28440 (76.6), maa://31400 (100.0), *maa://28650 (71.43)"
$ws.Range("P29").Value = "*maa://23168 (54.72), *maa://30050 (51.72)"
$ws.Range("AF29").Value = "*maa://24080 (69.05), ***maa://34960 (8.33), *maa://42865 (76.47)"
$ws.Range("AB30").Value = "maa://42979 (96.36), maa://45045 (100.0)"
$ws.Range("L31").Value = "maa://35926 (93.63), maa://36258 (83.33), *maa://43904 (77.78)"
$ws.Range("H32").Value = "maa://21895 (97.28), maa://36667 (98.39), **maa://20793 (38.78), maa://22760 (100.0)"
$ws.Range("T32").Value = "maa://42859 (96.34), maa://41108 (87.76), maa://41238 (96.2)"
$ws.Range("T34").Value = "maa://24526 (93.57)"
$ws.Range("L35").Value = "maa://41296 (96.67)"
$ws.Range("P37").Value = "maa://21280 (89.11), *maa://21239 (66.67)"
$ws.Range("P38").Value = "*maa://24383 (68.04)"
$ws.Range("H39").Value = "maa://25199 (84.82), maa://36670 (87.5), maa://30434 (89.39), ***maa://25036 (16.0), *maa://44165 (66.67), maa://45059 (100.0)"
$ws.Range("P39").Value = "maa://24709 (91.41)"
$ws.Range("H44").Value = "maa://29768 (97.83), maa://27728 (96.0)"
$ws.Range("H45").Value = "maa://21229 (84.86), maa://30807 (95.45), *maa://22767 (55.0), ***maa://20796 (13.79), *maa://42459 (77.78)"
$ws.Range("T45").Value = "**maa://39364 (38.1)"
$ws.Range("H55").Value = "maa://32532 (92.25)"
$ws.Range("H59").Value = "maa://27746 (83.02), maa://31270 (94.83)"
$ws.Range("H60").Value = "*maa://40438 (60.87)"
$ws.Range("H62").Value = "maa://42981 (96.43), maa://43903 (100.0)"
$ws.Range("H64").Value = "maa://44405 (91.3)"

# The following two cells hold plain numeric-looking text ("2" and "6")
# in the source data (style General, stored as text). A direct .Value
# assignment would auto-convert them to real numbers, so we force them
# back to text by setting NumberFormat to Text, assigning the value,
# then restoring the original General cell formatting by pasting the
# format from an unrelated General-formatted cell (C2) that uses the
# same original style index, avoiding creation of a brand-new style.
$ws.Range("AA30").NumberFormat = "@"
$ws.Range("AA30").Value = "2"
$ws.Range("C2").Copy() | Out-Null
$ws.Range("AA30").PasteSpecial(-4122) | Out-Null
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = "6"
$ws.Range("C2").Copy() | Out-Null
$ws.Range("G39").PasteSpecial(-4122) | Out-Null
